# Refresh the cryptos price/volume snapshot (GitHub Actions scheduled update).
# Source: coinranking.com snapshot refresh -- updates Price (D) and Volume(1h) (E)
# columns per coin row, and for two row pairs the ranking order of two coins
# swapped places so their full row (Coin/Link/Price/Volume) moved.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "70.587.52"
$ws.Range("E2").Value = "  +2.38%  "
# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.541.28"
$ws.Range("E3").Value = "  +1.44%  "
# Row 4
$ws.Range("E4").Value = "  -0.02%  "
# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "609.43"
$ws.Range("E5").Value = "  +5.10%  "
# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "173.34"
$ws.Range("E6").Value = "  +0.65%  "
# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.618"
$ws.Range("E7").Value = "  +0.75%  "
# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "3.540.20"
$ws.Range("E8").Value = "  +1.68%  "
# Row 9
$ws.Range("E9").Value = "  -0.08%  "
# Row 10
$ws.Range("E10").Value = "  +6.84%  "
# Row 11
$ws.Range("E11").Value = "  +1.62%  "
# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.586"
$ws.Range("E12").Value = "  -0.80%  "
# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "47.40"
$ws.Range("E13").Value = "  +2.05%  "
# Row 14
$ws.Range("E14").Value = "  +2.51%  "
# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "4.102.70"
$ws.Range("E15").Value = "  +1.19%  "
# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "627.44"
$ws.Range("E16").Value = "  -6.62%  "
# Row 17
$ws.Range("E17").Value = "  -2.53%  "
# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "70.601.61"
$ws.Range("E18").Value = "  +2.41%  "
# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "3.541.89"
$ws.Range("E19").Value = "  +1.28%  "
# Row 20
$ws.Range("E20").Value = "  -1.61%  "
# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "17.42"
$ws.Range("E21").Value = "  +0.71%  "
# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.886"
$ws.Range("E22").Value = "  -0.95%  "
# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "9.97"
$ws.Range("E23").Value = "  -10.16%  "
# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "15.93"
$ws.Range("E24").Value = "  -0.43%  "
# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "96.92"
$ws.Range("E25").Value = "  -0.31%  "
# Row 26
$ws.Range("E26").Value = "  +0.28%  "
# Row 27
$ws.Range("E27").Value = "  +0.05%  "
# Row 28
$ws.Range("E28").Value = "  -0.62%  "
# Row 29
$ws.Range("E29").Value = "  -1.17%  "
# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "33.52"
$ws.Range("E30").Value = "  +2.51%  "
# Row 31
$ws.Range("E31").Value = "  -1.59%  "
# Row 32
$ws.Range("E32").Value = "  -1.89%  "
# Row 33
$ws.Range("E33").Value = "  -1.17%  "
# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "7.01"
$ws.Range("E34").Value = "  -2.61%  "
# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "568.98"
$ws.Range("E35").Value = "  -4.37%  "
# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "3.61"
$ws.Range("E36").Value = "  +1.77%  "
# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "10.79"
$ws.Range("E37").Value = "  -0.16%  "
# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "57.72"
$ws.Range("E38").Value = "  +1.49%  "
# Row 39
$ws.Range("E39").Value = "  -1.19%  "
# Row 40
$ws.Range("E40").Value = "  -0.01%  "
# Row 41
$ws.Range("B41").Value = "VeChain"
$ws.Range("C41").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.0458"
$ws.Range("E41").Value = "  +5.41%  "
# Row 42
$ws.Range("B42").Value = "Kaspa"
$ws.Range("C42").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.143"
$ws.Range("E42").Value = "  +6.05%  "
# Row 43
$ws.Range("E43").Value = "  -1.15%  "
# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "3.341.49"
$ws.Range("E44").Value = "  -1.54%  "
# Row 45
$ws.Range("E45").Value = "  +5.55%  "
# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0₃0714"
$ws.Range("E46").Value = "  +1.95%  "
# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "33.10"
$ws.Range("E47").Value = "  +0.19%  "
# Row 48
$ws.Range("E48").Value = "  +3.25%  "
# Row 49
$ws.Range("E49").Value = "  -1.77%  "
# Row 50
$ws.Range("B50").Value = "Monero"
$ws.Range("C50").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "133.64"
$ws.Range("E50").Value = "  +0.56%  "
# Row 51
$ws.Range("B51").Value = "MXToken"
$ws.Range("C51").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "5.74"
$ws.Range("E51").Value = "  +0.19%  "
